$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"24.35712233333334"
$ws.Range("H2").Value = [double]"73.07136700000001"
$ws.Range("I2").Value = [double]"0.3750500562097488"
$ws.Range("J2").Value = [double]"0.3750500562097488"
$ws.Range("M2").Value = [double]"181.3526613333333"
$ws.Range("N2").Value = [double]"544.057984"
$ws.Range("O2").Value = [double]"0.9845849379007657"
$ws.Range("P2").Value = [double]"0.984584937900766"
$ws.Range("Q2").Value = [double]"4417.22895757157"
$ws.Range("R2").Value = [double]"39755.06061814413"
$ws.Range("S2").Value = [double]"0.3692686363029543"
$ws.Range("T2").Value = [double]"0.3692686363029543"

# Row 3
$ws.Range("G3").Value = [double]"24.35712233333334"
$ws.Range("H3").Value = [double]"73.07136700000001"
$ws.Range("I3").Value = [double]"0.3750500562097488"
$ws.Range("J3").Value = [double]"0.3750500562097488"
$ws.Range("O3").Value = [double]"0.003686045149950483"
$ws.Range("P3").Value = [double]"0.003686045149950484"
$ws.Range("Q3").Value = [double]"16.53702463699334"
$ws.Range("R3").Value = [double]"148.83322173294"
$ws.Range("S3").Value = [double]"0.001382451440680601"
$ws.Range("T3").Value = [double]"0.001382451440680601"

# Row 4
$ws.Range("G4").Value = [double]"24.35712233333334"
$ws.Range("H4").Value = [double]"73.07136700000001"
$ws.Range("I4").Value = [double]"0.3750500562097488"
$ws.Range("J4").Value = [double]"0.3750500562097488"
$ws.Range("M4").Value = [double]"0.6398506666666667"
$ws.Range("N4").Value = [double]"1.919552"
$ws.Range("O4").Value = [double]"0.003473824559694892"
$ws.Range("P4").Value = [double]"0.003473824559694892"
$ws.Range("Q4").Value = [double]"15.58492096306489"
$ws.Range("R4").Value = [double]"140.264288667584"
$ws.Range("S4").Value = [double]"0.001302858096376375"
$ws.Range("T4").Value = [double]"0.001302858096376375"

# Row 5
$ws.Range("G5").Value = [double]"24.35712233333334"
$ws.Range("H5").Value = [double]"73.07136700000001"
$ws.Range("I5").Value = [double]"0.3750500562097488"
$ws.Range("J5").Value = [double]"0.3750500562097488"
$ws.Range("M5").Value = [double]"1.520540333333333"
$ws.Range("N5").Value = [double]"4.561621"
$ws.Range("O5").Value = [double]"0.008255192389588805"
$ws.Range("P5").Value = [double]"0.008255192389588807"
$ws.Range("Q5").Value = [double]"37.03598691176744"
$ws.Range("R5").Value = [double]"333.323882205907"
$ws.Range("S5").Value = [double]"0.003096110369737573"
$ws.Range("T5").Value = [double]"0.003096110369737573"

# Row 6
$ws.Range("I6").Value = [double]"0.2805618708302703"
$ws.Range("J6").Value = [double]"0.2805618708302702"
$ws.Range("M6").Value = [double]"181.3526613333333"
$ws.Range("N6").Value = [double]"544.057984"
$ws.Range("O6").Value = [double]"0.9845849379007657"
$ws.Range("P6").Value = [double]"0.984584937900766"
$ws.Range("Q6").Value = [double]"3304.374975293525"
$ws.Range("R6").Value = [double]"29739.37477764173"
$ws.Range("S6").Value = [double]"0.2762369921687443"
$ws.Range("T6").Value = [double]"0.2762369921687443"

# Row 7
$ws.Range("I7").Value = [double]"0.2805618708302703"
$ws.Range("J7").Value = [double]"0.2805618708302702"
$ws.Range("O7").Value = [double]"0.003686045149950483"
$ws.Range("P7").Value = [double]"0.003686045149950484"
$ws.Range("S7").Value = [double]"0.001034163723234952"
$ws.Range("T7").Value = [double]"0.001034163723234952"

# Row 8
$ws.Range("I8").Value = [double]"0.2805618708302703"
$ws.Range("J8").Value = [double]"0.2805618708302702"
$ws.Range("M8").Value = [double]"0.6398506666666667"
$ws.Range("N8").Value = [double]"1.919552"
$ws.Range("O8").Value = [double]"0.003473824559694892"
$ws.Range("P8").Value = [double]"0.003473824559694892"
$ws.Range("Q8").Value = [double]"11.65853600004267"
$ws.Range("R8").Value = [double]"104.926824000384"
$ws.Range("S8").Value = [double]"0.0009746227174041388"
$ws.Range("T8").Value = [double]"0.0009746227174041388"

# Row 9
$ws.Range("I9").Value = [double]"0.2805618708302703"
$ws.Range("J9").Value = [double]"0.2805618708302702"
$ws.Range("M9").Value = [double]"1.520540333333333"
$ws.Range("N9").Value = [double]"4.561621"
$ws.Range("O9").Value = [double]"0.008255192389588805"
$ws.Range("P9").Value = [double]"0.008255192389588807"
$ws.Range("Q9").Value = [double]"27.70533053913133"
$ws.Range("R9").Value = [double]"249.347974852182"
$ws.Range("S9").Value = [double]"0.002316092220886845"
$ws.Range("T9").Value = [double]"0.002316092220886845"

# Row 10
$ws.Range("G10").Value = [double]"22.31748066666667"
$ws.Range("H10").Value = [double]"66.952442"
$ws.Range("I10").Value = [double]"0.3436437303202491"
$ws.Range("J10").Value = [double]"0.343643730320249"
$ws.Range("M10").Value = [double]"181.3526613333333"
$ws.Range("N10").Value = [double]"544.057984"
$ws.Range("O10").Value = [double]"0.9845849379007657"
$ws.Range("P10").Value = [double]"0.984584937900766"
$ws.Range("Q10").Value = [double]"4047.334513155215"
$ws.Range("R10").Value = [double]"36426.01061839693"
$ws.Range("S10").Value = [double]"0.3383464408773499"
$ws.Range("T10").Value = [double]"0.3383464408773499"

# Row 11
$ws.Range("G11").Value = [double]"22.31748066666667"
$ws.Range("H11").Value = [double]"66.952442"
$ws.Range("I11").Value = [double]"0.3436437303202491"
$ws.Range("J11").Value = [double]"0.343643730320249"
$ws.Range("O11").Value = [double]"0.003686045149950483"
$ws.Range("P11").Value = [double]"0.003686045149950484"
$ws.Range("Q11").Value = [double]"15.15223032382667"
$ws.Range("R11").Value = [double]"136.37007291444"
$ws.Range("S11").Value = [double]"0.001266686305457846"
$ws.Range("T11").Value = [double]"0.001266686305457846"

# Row 12
$ws.Range("G12").Value = [double]"22.31748066666667"
$ws.Range("H12").Value = [double]"66.952442"
$ws.Range("I12").Value = [double]"0.3436437303202491"
$ws.Range("J12").Value = [double]"0.343643730320249"
$ws.Range("M12").Value = [double]"0.6398506666666667"
$ws.Range("N12").Value = [double]"1.919552"
$ws.Range("O12").Value = [double]"0.003473824559694892"
$ws.Range("P12").Value = [double]"0.003473824559694892"
$ws.Range("Q12").Value = [double]"14.27985488288711"
$ws.Range("R12").Value = [double]"128.518693945984"
$ws.Range("S12").Value = [double]"0.00119375803017165"
$ws.Range("T12").Value = [double]"0.001193758030171649"

# Row 13
$ws.Range("G13").Value = [double]"22.31748066666667"
$ws.Range("H13").Value = [double]"66.952442"
$ws.Range("I13").Value = [double]"0.3436437303202491"
$ws.Range("J13").Value = [double]"0.343643730320249"
$ws.Range("M13").Value = [double]"1.520540333333333"
$ws.Range("N13").Value = [double]"4.561621"
$ws.Range("O13").Value = [double]"0.008255192389588805"
$ws.Range("P13").Value = [double]"0.008255192389588807"
$ws.Range("Q13").Value = [double]"33.93462949205355"
$ws.Range("R13").Value = [double]"305.411665428482"
$ws.Range("S13").Value = [double]"0.002836845107269628"
$ws.Range("T13").Value = [double]"0.002836845107269628"

# Row 14
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.04834033333333334"
$ws.Range("H14").Value = [double]"0.145021"
$ws.Range("I14").Value = [double]"0.0007443426397318391"
$ws.Range("J14").Value = [double]"0.0007443426397318388"
$ws.Range("M14").Value = [double]"181.3526613333333"
$ws.Range("N14").Value = [double]"544.057984"
$ws.Range("O14").Value = [double]"0.9845849379007657"
$ws.Range("P14").Value = [double]"0.984584937900766"
$ws.Range("Q14").Value = [double]"8.766648099740445"
$ws.Range("R14").Value = [double]"78.89983289766401"
$ws.Range("S14").Value = [double]"0.0007328685517172648"
$ws.Range("T14").Value = [double]"0.0007328685517172648"

# Row 15
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.04834033333333334"
$ws.Range("H15").Value = [double]"0.145021"
$ws.Range("I15").Value = [double]"0.0007443426397318391"
$ws.Range("J15").Value = [double]"0.0007443426397318388"
$ws.Range("O15").Value = [double]"0.003686045149950483"
$ws.Range("P15").Value = [double]"0.003686045149950484"
$ws.Range("Q15").Value = [double]"0.03282018591333333"
$ws.Range("R15").Value = [double]"0.29538167322"
$ws.Range("S15").Value = [double]"2.743680577084885E-06"
$ws.Range("T15").Value = [double]"2.743680577084885E-06"

# Row 16
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.04834033333333334"
$ws.Range("H16").Value = [double]"0.145021"
$ws.Range("I16").Value = [double]"0.0007443426397318391"
$ws.Range("J16").Value = [double]"0.0007443426397318388"
$ws.Range("M16").Value = [double]"0.6398506666666667"
$ws.Range("N16").Value = [double]"1.919552"
$ws.Range("O16").Value = [double]"0.003473824559694892"
$ws.Range("P16").Value = [double]"0.003473824559694892"
$ws.Range("Q16").Value = [double]"0.03093059451022223"
$ws.Range("R16").Value = [double]"0.278375350592"
$ws.Range("S16").Value = [double]"2.58571574272859E-06"
$ws.Range("T16").Value = [double]"2.585715742728589E-06"

# Row 17
$ws.Range("E17").Value = [double]"2"
$ws.Range("F17").Value = [double]"0.6666666666666666"
$ws.Range("G17").Value = [double]"0.04834033333333334"
$ws.Range("H17").Value = [double]"0.145021"
$ws.Range("I17").Value = [double]"0.0007443426397318391"
$ws.Range("J17").Value = [double]"0.0007443426397318388"
$ws.Range("M17").Value = [double]"1.520540333333333"
$ws.Range("N17").Value = [double]"4.561621"
$ws.Range("O17").Value = [double]"0.008255192389588805"
$ws.Range("P17").Value = [double]"0.008255192389588807"
$ws.Range("Q17").Value = [double]"0.07350342656011111"
$ws.Range("R17").Value = [double]"0.661530839041"
$ws.Range("S17").Value = [double]"6.14469169476072E-06"
$ws.Range("T17").Value = [double]"6.144691694760719E-06"
